$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '65.305.58'
Set-TextValue $ws.Range('E2') '  -1.95%  '
Set-TextValue $ws.Range('D3') '3.508.79'
Set-TextValue $ws.Range('E3') '  -2.13%  '
Set-TextValue $ws.Range('E4') '  -0.06%  '
Set-TextValue $ws.Range('D5') '599.71'
Set-TextValue $ws.Range('E5') '  -1.59%  '
Set-TextValue $ws.Range('D6') '142.70'
Set-TextValue $ws.Range('E6') '  -2.95%  '
Set-TextValue $ws.Range('D7') '3.506.63'
Set-TextValue $ws.Range('E7') '  -2.14%  '
Set-TextValue $ws.Range('E8') '  -0.16%  '
Set-TextValue $ws.Range('D9') '0.517'
Set-TextValue $ws.Range('E9') '  +5.38%  '
Set-TextValue $ws.Range('E10') '  -2.85%  '
Set-TextValue $ws.Range('D11') '7.82'
Set-TextValue $ws.Range('E11') '  -2.36%  '
Set-TextValue $ws.Range('D12') '0.402'
Set-TextValue $ws.Range('E12') '  -3.07%  '
Set-TextValue $ws.Range('D13') '4.120.75'
Set-TextValue $ws.Range('E13') '  -1.75%  '
Set-TextValue $ws.Range('D14') '0.0000196'
Set-TextValue $ws.Range('E14') '  -6.37%  '
Set-TextValue $ws.Range('D15') '28.32'
Set-TextValue $ws.Range('E15') '  -5.62%  '
Set-TextValue $ws.Range('D16') '3.513.02'
Set-TextValue $ws.Range('E16') '  -2.04%  '
Set-TextValue $ws.Range('E17') '  +1.35%  '
Set-TextValue $ws.Range('D18') '65.382.22'
Set-TextValue $ws.Range('E18') '  -1.97%  '
Set-TextValue $ws.Range('D19') '10.89'
Set-TextValue $ws.Range('E19') '  -4.40%  '
Set-TextValue $ws.Range('D20') '6.14'
Set-TextValue $ws.Range('E20') '  -2.30%  '
Set-TextValue $ws.Range('D21') '14.35'
Set-TextValue $ws.Range('E21') '  -4.89%  '
Set-TextValue $ws.Range('D22') '416.82'
Set-TextValue $ws.Range('E22') '  -3.60%  '
Set-TextValue $ws.Range('D23') '0.594'
Set-TextValue $ws.Range('E23') '  -4.46%  '
Set-TextValue $ws.Range('D24') '76.94'
Set-TextValue $ws.Range('E24') '  -2.68%  '
Set-TextValue $ws.Range('D25') '3.654.58'
Set-TextValue $ws.Range('E25') '  -2.12%  '
Set-TextValue $ws.Range('E26') '  -0.01%  '
Set-TextValue $ws.Range('D27') '0.0000113'
Set-TextValue $ws.Range('E27') '  -6.13%  '
Set-TextValue $ws.Range('D28') '2.44'
Set-TextValue $ws.Range('E28') '  -2.89%  '
Set-TextValue $ws.Range('D29') '7.71'
Set-TextValue $ws.Range('E29') '  -5.52%  '
Set-TextValue $ws.Range('D30') '8.82'
Set-TextValue $ws.Range('E30') '  -5.14%  '
Set-TextValue $ws.Range('D31') '1.00'
Set-TextValue $ws.Range('E31') '  +0.02%  '
Set-TextValue $ws.Range('D32') '3.519.16'
Set-TextValue $ws.Range('E32') '  -1.73%  '
Set-TextValue $ws.Range('D33') '0.154'
Set-TextValue $ws.Range('E33') '  -0.94%  '
Set-TextValue $ws.Range('D34') '24.09'
Set-TextValue $ws.Range('E34') '  -5.65%  '
Set-TextValue $ws.Range('D36') '1.31'
Set-TextValue $ws.Range('E36') '  -9.44%  '
Set-TextValue $ws.Range('D37') '7.49'
Set-TextValue $ws.Range('E37') '  -4.48%  '
Set-TextValue $ws.Range('B38') 'Monero'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D38') '172.02'
Set-TextValue $ws.Range('E38') '  -1.03%  '
Set-TextValue $ws.Range('B39') 'NEARProtocol'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D39') '5.20'
Set-TextValue $ws.Range('E39') '  -7.40%  '
Set-TextValue $ws.Range('E40') '  -8.72%  '
Set-TextValue $ws.Range('D41') '0.0808'
Set-TextValue $ws.Range('E41') '  -5.48%  '
Set-TextValue $ws.Range('D42') '5.00'
Set-TextValue $ws.Range('E42') '  -4.35%  '
Set-TextValue $ws.Range('D43') '0.853'
Set-TextValue $ws.Range('E43') '  -4.73%  '
Set-TextValue $ws.Range('D44') '45.14'
Set-TextValue $ws.Range('E44') '  -2.04%  '
Set-TextValue $ws.Range('D45') '1.76'
Set-TextValue $ws.Range('E45') '  -8.03%  '
Set-TextValue $ws.Range('E46') '  +0.09%  '
Set-TextValue $ws.Range('B47') 'dogwifhat'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D47') '2.34'
Set-TextValue $ws.Range('E47') '  -8.46%  '
Set-TextValue $ws.Range('B48') 'EnergySwap'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D48') '23.16'
Set-TextValue $ws.Range('E48') '  -2.14%  '
Set-TextValue $ws.Range('D49') '7.00'
Set-TextValue $ws.Range('E49') '  -2.89%  '
Set-TextValue $ws.Range('D50') '1.09'
Set-TextValue $ws.Range('E50') '  -8.43%  '
Set-TextValue $ws.Range('D51') '0.898'
Set-TextValue $ws.Range('E51') '  -5.12%  '
